$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of test case data appended to the "Test Cases" sheet.
$tcids = @("RCC001","RCC002","RCC009","RCC014","RCC015","RCC016","RCC017","RCC018","RCC019","RCC021","RCC022","RCC023","RCC024")

$startRow = 30
for ($i = 0; $i -lt $tcids.Count; $i++) {
  $r = $startRow + $i
  # Copy the plain formatting (style) used by row 2 onto the new row.
  $ws.Range("A2:E2").Copy()
  $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)

  $ws.Cells.Item($r, 1).Value = $tcids[$i]
  $ws.Cells.Item($r, 2).Value = "TBD"
  $ws.Cells.Item($r, 3).Value = "TBD"
  $ws.Cells.Item($r, 4).Value = "Y"
}

# Update the active selection on the sheet to C16.
$ws.Range("C16").Select() | Out-Null
